$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.96675092036255
$ws.Range("D2").Value = 7.789308099070766
$ws.Range("H2").Value = 6.486603965511964
$ws.Range("J2").Value = 5.557842537230018
